# Final pre-jettison cleanup of the NUFORC sightings sheet:
#  - Column C ("Time") was stored as free-text strings (e.g. "19:05:00")
#    that didn't actually correspond to the sighting rows; it is replaced
#    with the correct numeric time-of-day serials, formatted h:mm:ss.
#    (This also lets the now-unused text-time shared strings fall out of
#    the shared string table on save.)
#  - Column B ("Date") gets a best-fit width now that column C no longer
#    forces an artificially wide column.
#  - Leave the selection on C16, matching where editing finished.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timeValues = @(
    "4.1666666666666664E-2","8.3333333333333329E-2","0.125","0.16666666666666666","0.20833333333333401",
    "0.25","0.29166666666666702","0.33333333333333398","0.375","0.41666666666666702",
    "0.45833333333333398","0.5","0.54166666666666696","0.54166666666666663","0.625",
    "0.66666666666666696","0.70833333333333404","0.75","0.79166666666666696","0.83333333333333404",
    "0.875","0.91666666666666696","0.95833333333333404","1","1.0416666666666701",
    "1.0833333333333399","0.54166666666666663","1.1666666666666701","1.2083333333333399","1.25",
    "0.54166666666666663","1.3333333333333399","1.375","1.4166666666666701","1.4583333333333399",
    "1.5","1.5416666666666701"
)

$firstRow = 2
for ($i = 0; $i -lt $timeValues.Length; $i++) {
    $row = $firstRow + $i
    $cell = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "h:mm:ss"
    $cell.Value = [double]$timeValues[$i]
}

# Best-fit the Date column now that its content governs the visible width.
$ws.Columns.Item(2).AutoFit()

# Restore the selection left active at the end of the edit session.
$ws.Range("C16").Select()
